# IK construction: update IK-related rows/cells in Actor_RobotObj_01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 4..9 : W column (IK blend?) 0.5 -> 0, AC column (deg) 45 -> 0 ---
foreach ($r in 4..9) {
    $ws.Cells.Item($r, 23).Value = 0   # column W
    $ws.Cells.Item($r, 29).Value = 0   # column AC
}

# --- Row 11 : new IK end-effector object entry (objEE01 / obj6, " IK object") ---
$ws.Range("A11").Value = 100
$ws.Range("B11").Value = "objEE01"
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = "obj6"
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 200
$ws.Range("W11").Value = 0
$ws.Range("X11").Value = 0
$ws.Range("Y11").Value = 0
$ws.Range("Z11").Value = 1
$ws.Range("AA11").Value = 0
$ws.Range("AB11").Value = 0
$ws.Range("AC11").Value = 0
$ws.Range("AD11").Value = 1
$ws.Range("AE11").Value = " IK object"

# --- Row 13 : new IK target entry (target01 / root, "IK target") ---
$ws.Range("A13").Value = 200
$ws.Range("B13").Value = "target01"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = "root"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 1
$ws.Range("U13").Value = 200
$ws.Range("W13").Value = 0
$ws.Range("X13").Value = 0
$ws.Range("Y13").Value = 0
$ws.Range("Z13").Value = 1
$ws.Range("AA13").Value = 0
$ws.Range("AB13").Value = 0
$ws.Range("AC13").Value = 0
$ws.Range("AD13").Value = 1
$ws.Range("AE13").Value = "IK target"

# --- sheet view: selection moved from W4:W9 to AC4:AC9 ---
$ws.Range("AC4:AC9").Select()

# --- workbook window geometry ---
$excel.Windows.Item(1).Left = -23535
$excel.Windows.Item(1).Top = 4905
$excel.Windows.Item(1).Width = 17220
$excel.Windows.Item(1).Height = 10980
